$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "boson" (col E), shifting columns E:N
# (boson..col) one place to the right, to F:O.
$ws.Columns("E:E").Insert()

# Populate the newly inserted "pt_max" column: header + a constant 50
# for every data row.
$ws.Range("E1").Value = "pt_max"
$ws.Range("E2:E17").Value = 50

# The old "syst1_c" header (now shifted to column I) is renamed to
# "syst1_u"; the underlying numeric data in that column is unchanged.
$ws.Range("I1").Value = "syst1_u"

# Leave the selection where the author ended up after editing.
[void]$ws.Range("E22").Select()
